$wb = $excel.ActiveWorkbook
Write-Host $wb.ActiveSheet.Name
